$wb = $excel.ActiveWorkbook

# --- Rename header cells on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after "Monthly Trend" (i.e. at the end) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match the bold/bordered header formatting used on the other sheets
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$data = @()
$data += ,@(45333.99999999999, 0, -1569.172289935177, 752.8520966964285)
$data += ,@(45340.99999999999, 0, -1394.342137087983, 716.7658629600322)
$data += ,@(45347.99999999999, 0, -1408.493556663553, 749.1275843261839)
$data += ,@(45354.99999999999, 0, -1271.673651289975, 875.6580862597027)
$data += ,@(45361.99999999999, 0, -1257.853607410123, 987.152481912521)
$data += ,@(45375.99999999999, 0, -1049.92043266749, 1145.695957807578)
$data += ,@(45382.99999999999, 73, -1030.620185438274, 1111.184468174231)
$data += ,@(45389.99999999999, 145, -887.9137387637942, 1285.43044318444)
$data += ,@(45396.99999999999, 218, -937.7054572448545, 1254.269491717867)
$data += ,@(45403.99999999999, 291, -829.9028129368346, 1373.321394635203)
$data += ,@(45410.99999999999, 363, -825.4883685708548, 1381.31428395441)
$data += ,@(45417.99999999999, 436, -614.0397482861968, 1504.805525162318)
$data += ,@(45424.99999999999, 508, -631.2064384562534, 1603.295794895417)
$data += ,@(45438.99999999999, 654, -436.0178032084112, 1770.030609216602)
$data += ,@(45445.99999999999, 726, -391.9108860178037, 1815.42405780285)
$data += ,@(45452.99999999999, 799, -280.7714665989869, 1876.553498237407)
$data += ,@(45459.99999999999, 872, -256.8469479529409, 1953.456535558113)
$data += ,@(45466.99999999999, 944, -196.5461739519828, 1998.761153346092)
$data += ,@(45473.99999999999, 1017, -74.17804289807744, 2109.539042813391)
$data += ,@(45480.99999999999, 1089, 19.91974239474641, 2213.528522994807)
$data += ,@(45487.99999999999, 1162, 68.82891172958142, 2286.216103236714)
$data += ,@(45501.99999999999, 1307, 205.6806606252317, 2391.206511018382)
$data += ,@(45508.99999999999, 1380, 286.4663823180388, 2465.81814572166)
$data += ,@(45522.99999999999, 1525, 351.2823458870184, 2677.038398988386)
$data += ,@(45536.99999999999, 1670, 572.1952836082592, 2806.062084353113)
$data += ,@(45543.99999999999, 1743, 642.7468034647292, 2803.540887852565)
$data += ,@(45550.99999999999, 1815, 736.2002679944483, 2902.715612731477)
$data += ,@(45557.99999999999, 1888, 906.20972570149, 3047.607477291248)
$data += ,@(45564.99999999999, 1961, 808.7646297415231, 3044.237853110549)
$data += ,@(45571.99999999999, 2033, 968.8392044025709, 3116.044054770035)
$data += ,@(45585.99999999999, 2179, 1105.680837672452, 3326.149952662897)
$data += ,@(45592.99999999999, 2251, 1128.847116340645, 3379.577327465195)
$data += ,@(45599.99999999999, 2324, 1222.875146798304, 3455.04242110629)
$data += ,@(45606.99999999999, 2396, 1362.980230838848, 3462.218718208154)
$data += ,@(45613.99999999999, 2469, 1298.737483617479, 3564.740738077106)
$data += ,@(45620.99999999999, 2542, 1330.0635677506, 3616.636936475835)
$data += ,@(45627.99999999999, 2614, 1544.049182224908, 3706.339389376767)
$data += ,@(45634.99999999999, 2687, 1522.207740797258, 3715.096052354489)
$data += ,@(45641.99999999999, 2759, 1627.875940224645, 3827.018805182168)
$data += ,@(45648.99999999999, 2832, 1763.062832029964, 3877.807310066887)
$data += ,@(45655.99999999999, 2905, 1821.712842408888, 4046.244620727868)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Match the date-number-format used for column A on the other sheets
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A42").PasteSpecial(-4122)

$excel.CutCopyMode = $false
